# Auto-generated edit script applying numeric cell updates per commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J29").Value = 4003
$ws.Range("H29").Value = 166668060
$ws.Range("L29").Value = 12009
$ws.Range("N29").Value = -12571
$ws.Range("I38").Value = 833.73334
$ws.Range("K38").Value = 2501.20002
$ws.Range("M38").Value = -2129.20002
$ws.Range("H38").Value = 2931.9546
$ws.Range("N42").Value = -2434
$ws.Range("J42").Value = 658
$ws.Range("M42").Value = -200000242
$ws.Range("K42").Value = 200000472
$ws.Range("I42").Value = 66666824
$ws.Range("L42").Value = 1974
$ws.Range("H42").Value = 55555796
$ws.Range("M51").Value = -1014.4286
$ws.Range("I51").Value = 1498.4286
$ws.Range("K51").Value = 1498.4286
$ws.Range("H51").Value = 1999
$ws.Range("M132").Value = -10595
$ws.Range("I132").Value = 4375
$ws.Range("H132").Value = 3899.5
$ws.Range("J132").Value = 3582.5
$ws.Range("L132").Value = 10747.5
$ws.Range("K132").Value = 13125
$ws.Range("N132").Value = -15807.5
$ws.Range("K137").Value = 3925.875
$ws.Range("I137").Value = 1308.625
$ws.Range("M137").Value = -1375.875
$ws.Range("H137").Value = 1481.3182
$ws.Range("L138").Value = 9542.206200000001
$ws.Range("I138").Value = 1470.0526
$ws.Range("J138").Value = 3180.7354
$ws.Range("M138").Value = 729.8422
$ws.Range("N138").Value = -19822.2062
$ws.Range("K138").Value = 4410.1578
$ws.Range("H138").Value = 2567.4717
$ws.Range("N139").Value = -85280
$ws.Range("L139").Value = 75000
$ws.Range("J139").Value = 75000
$ws.Range("H139").Value = 75000
$ws.Range("M140").ClearContents()
$ws.Range("I140").Value = 0
$ws.Range("H140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("J140").Value = 0
$ws.Range("I141").Value = 3684.3076
$ws.Range("K141").Value = 11052.9228
$ws.Range("H141").Value = 3962.1875
$ws.Range("M141").Value = -5872.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N2").Value = -4479.25
$ws.Range("M2").Value = -1877
$ws.Range("J2").Value = 4253.25
$ws.Range("I2").Value = 1990
$ws.Range("K2").Value = 1990
$ws.Range("H2").Value = 3121.625
$ws.Range("L2").Value = 4253.25
$ws.Range("K32").Value = 2928.9473
$ws.Range("M32").Value = -2641.9473
$ws.Range("I32").Value = 2928.9473
$ws.Range("H32").Value = 3063.5557
$ws.Range("N116").Value = -8841.25
$ws.Range("H116").Value = 3121.625
$ws.Range("I116").Value = 1990
$ws.Range("J116").Value = 4253.25
$ws.Range("L116").Value = 4253.25
$ws.Range("K116").Value = 1990
$ws.Range("M116").Value = 304
$ws.Range("M132").Value = -4182.6362
$ws.Range("I132").Value = 2237.5454
$ws.Range("H132").Value = 2237.5454
$ws.Range("K132").Value = 6712.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J3").Value = 4253.25
$ws.Range("M3").Value = -1876
$ws.Range("I3").Value = 1990
$ws.Range("L3").Value = 4253.25
$ws.Range("N3").Value = -4481.25
$ws.Range("K3").Value = 1990
$ws.Range("H3").Value = 3121.625
$ws.Range("K20").Value = 852.5
$ws.Range("H20").Value = 1562.4615
$ws.Range("I20").Value = 852.5
$ws.Range("M20").Value = -605.5
$ws.Range("H134").Value = 1906.5333
$ws.Range("I134").Value = 1618.091
$ws.Range("M134").Value = -2319.272999999999
$ws.Range("K134").Value = 4854.272999999999
$ws.Range("N135").Value = -96138
$ws.Range("L135").Value = 85998
$ws.Range("H135").Value = 85998
$ws.Range("J135").Value = 85998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M29").Value = -1806
$ws.Range("I29").Value = 2099
$ws.Range("K29").Value = 2099
$ws.Range("J29").Value = 0
$ws.Range("H29").Value = 2099
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("L51").Value = 40000
$ws.Range("N51").Value = -41472
$ws.Range("J51").Value = 40000
$ws.Range("H51").Value = 79316.664
$ws.Range("N58").Value = -1950
$ws.Range("I58").Value = 1615.6875
$ws.Range("L58").Value = 1544
$ws.Range("K58").Value = 1615.6875
$ws.Range("H58").Value = 1586.4814
$ws.Range("M58").Value = -1412.6875
$ws.Range("J58").Value = 1544
$ws.Range("L61").Value = 40000
$ws.Range("N61").Value = -40696
$ws.Range("H61").Value = 79316.664
$ws.Range("J61").Value = 40000
$ws.Range("J134").Value = 2250
$ws.Range("H134").Value = 3153.0908
$ws.Range("N134").Value = -11820
$ws.Range("L134").Value = 6750
$ws.Range("K136").Value = 4847.0625
$ws.Range("J136").Value = 1544
$ws.Range("I136").Value = 1615.6875
$ws.Range("L136").Value = 4632
$ws.Range("H136").Value = 1586.4814
$ws.Range("N136").Value = -9732
$ws.Range("M136").Value = -2297.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K8").Value = 119997
$ws.Range("I8").Value = 39999
$ws.Range("H8").Value = 39999
$ws.Range("M8").Value = -119858
$ws.Range("M12").Value = 7.25
$ws.Range("J12").Value = 10
$ws.Range("K12").Value = 165.75
$ws.Range("L12").Value = 30
$ws.Range("N12").Value = -376
$ws.Range("I12").Value = 55.25
$ws.Range("H12").Value = 46.2
$ws.Range("H39").Value = 149.33333
$ws.Range("N39").Value = -1038
$ws.Range("L39").Value = 450
$ws.Range("J39").Value = 150
$ws.Range("H125").Value = 9999
$ws.Range("I125").Value = 9999
$ws.Range("K125").Value = 29997
$ws.Range("M125").Value = -25077
$ws.Range("M131").Value = 2940.9999
$ws.Range("H131").Value = 1631.8846
$ws.Range("N131").Value = -16456.2357
$ws.Range("J131").Value = 2125.4119
$ws.Range("I131").Value = 699.6667
$ws.Range("K131").Value = 2099.0001
$ws.Range("L131").Value = 6376.2357

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M70").Value = -15854.875
$ws.Range("H70").Value = 13899.8
$ws.Range("I70").Value = 16124.875
$ws.Range("K70").Value = 16124.875
$ws.Range("H73").Value = 13899.8
$ws.Range("I73").Value = 16124.875
$ws.Range("K73").Value = 16124.875
$ws.Range("M73").Value = -15188.875
$ws.Range("I97").Value = 700
$ws.Range("K97").Value = 700
$ws.Range("N97").Value = -1292
$ws.Range("M97").Value = -204
$ws.Range("J97").Value = 300
$ws.Range("H97").Value = 642.8570999999999
$ws.Range("L97").Value = 300
$ws.Range("M122").Value = -9045.000100000001
$ws.Range("H122").Value = 4776.5557
$ws.Range("J122").Value = 6666.3335
$ws.Range("N122").Value = -24899.0005
$ws.Range("K122").Value = 11495.0001
$ws.Range("L122").Value = 19999.0005
$ws.Range("I122").Value = 3831.6667
$ws.Range("M132").Value = -1986.071599999999
$ws.Range("I132").Value = 1505.3572
$ws.Range("H132").Value = 2424.3157
$ws.Range("J132").Value = 4997.4
$ws.Range("L132").Value = 14992.2
$ws.Range("K132").Value = 4516.071599999999
$ws.Range("N132").Value = -20052.2
$ws.Range("L141").Value = 80000
$ws.Range("N141").Value = -90360
$ws.Range("J141").Value = 80000
$ws.Range("H141").Value = 80000

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L46").Value = 1597.8334
$ws.Range("H46").Value = 1687
$ws.Range("M46").Value = -1677.3334
$ws.Range("I46").Value = 1865.3334
$ws.Range("K46").Value = 1865.3334
$ws.Range("N46").Value = -1973.8334
$ws.Range("J46").Value = 1597.8334
$ws.Range("I100").Value = 5116.3335
$ws.Range("K100").Value = 5116.3335
$ws.Range("J100").Value = 7599.7144
$ws.Range("H100").Value = 6031.263
$ws.Range("L100").Value = 7599.7144
$ws.Range("M100").Value = -4575.3335
$ws.Range("N100").Value = -8681.714400000001
$ws.Range("M132").Value = -4569.5
$ws.Range("I132").Value = 2366.5
$ws.Range("H132").Value = 4683.25
$ws.Range("K132").Value = 7099.5
$ws.Range("K136").Value = 5183.5002
$ws.Range("I136").Value = 1727.8334
$ws.Range("H136").Value = 2926.5
$ws.Range("M136").Value = -2633.5002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M122").Value = -8434.375
$ws.Range("H122").Value = 6078.846
$ws.Range("J122").Value = 10000
$ws.Range("N122").Value = -34900
$ws.Range("K122").Value = 10884.375
$ws.Range("L122").Value = 30000
$ws.Range("I122").Value = 3628.125
$ws.Range("I126").Value = 1221.3334
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3664.0002
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1194.0002
$ws.Range("N126").ClearContents()
$ws.Range("H126").Value = 1221.3334
$ws.Range("K136").Value = 12601.2
$ws.Range("J136").Value = 4542.3335
$ws.Range("I136").Value = 4200.4
$ws.Range("L136").Value = 13627.0005
$ws.Range("H136").Value = 4386.909
$ws.Range("N136").Value = -18727.0005
$ws.Range("M136").Value = -10051.2
$ws.Range("L138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").ClearContents()
$ws.Range("K138").Value = 0
$ws.Range("H138").Value = 0
